$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B and C (constant for all data rows): Ligand symbol = Lama3, Receptor symbol = Sdc2
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = "Lama3"
    $ws.Cells.Item($r, 3).Value = "Sdc2"
}

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 5.624269666666667
$ws.Cells.Item(2, 8).Value = 16.872809
$ws.Cells.Item(2, 9).Value = 0.7456305897517299
$ws.Cells.Item(2, 10).Value = 0.74563058975173
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 1.321445333333333
$ws.Cells.Item(2, 14).Value = 3.964336
$ws.Cells.Item(2, 15).Value = 0.01021782062667047
$ws.Cells.Item(2, 16).Value = 0.01021782062667047
$ws.Cells.Item(2, 17).Value = 7.432164904424889
$ws.Cells.Item(2, 18).Value = 66.889484139824
$ws.Cells.Item(2, 19).Value = 0.007618719619841691
$ws.Cells.Item(2, 20).Value = 0.007618719619841694

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 5.624269666666667
$ws.Cells.Item(3, 8).Value = 16.872809
$ws.Cells.Item(3, 9).Value = 0.7456305897517299
$ws.Cells.Item(3, 10).Value = 0.74563058975173
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 105.9632263333333
$ws.Cells.Item(3, 14).Value = 317.889679
$ws.Cells.Item(3, 15).Value = 0.819340166699254
$ws.Cells.Item(3, 16).Value = 0.8193401666992541
$ws.Cells.Item(3, 17).Value = 595.9657596487012
$ws.Cells.Item(3, 18).Value = 5363.691836838311
$ws.Cells.Item(3, 19).Value = 0.6109250917032455
$ws.Cells.Item(3, 20).Value = 0.6109250917032456

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 5.624269666666667
$ws.Cells.Item(4, 8).Value = 16.872809
$ws.Cells.Item(4, 9).Value = 0.7456305897517299
$ws.Cells.Item(4, 10).Value = 0.74563058975173
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 22.04284166666666
$ws.Cells.Item(4, 14).Value = 66.128525
$ws.Cells.Item(4, 15).Value = 0.1704420126740755
$ws.Cells.Item(4, 16).Value = 0.1704420126740755
$ws.Cells.Item(4, 17).Value = 123.9748857529694
$ws.Cells.Item(4, 18).Value = 1115.773971776725
$ws.Cells.Item(4, 19).Value = 0.1270867784286427
$ws.Cells.Item(4, 20).Value = 0.1270867784286427

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.1400846666666667
$ws.Cells.Item(5, 8).Value = 0.420254
$ws.Cells.Item(5, 9).Value = 0.01857155129685422
$ws.Cells.Item(5, 10).Value = 0.01857155129685422
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 1.321445333333333
$ws.Cells.Item(5, 14).Value = 3.964336
$ws.Cells.Item(5, 15).Value = 0.01021782062667047
$ws.Cells.Item(5, 16).Value = 0.01021782062667047
$ws.Cells.Item(5, 17).Value = 0.1851142290382222
$ws.Cells.Item(5, 18).Value = 1.666028061344
$ws.Cells.Item(5, 19).Value = 0.0001897607799102657
$ws.Cells.Item(5, 20).Value = 0.0001897607799102658

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.1400846666666667
$ws.Cells.Item(6, 8).Value = 0.420254
$ws.Cells.Item(6, 9).Value = 0.01857155129685422
$ws.Cells.Item(6, 10).Value = 0.01857155129685422
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 105.9632263333333
$ws.Cells.Item(6, 14).Value = 317.889679
$ws.Cells.Item(6, 15).Value = 0.819340166699254
$ws.Cells.Item(6, 16).Value = 0.8193401666992541
$ws.Cells.Item(6, 17).Value = 14.84382323982956
$ws.Cells.Item(6, 18).Value = 133.594409158466
$ws.Cells.Item(6, 19).Value = 0.01521641793542828
$ws.Cells.Item(6, 20).Value = 0.01521641793542829

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.1400846666666667
$ws.Cells.Item(7, 8).Value = 0.420254
$ws.Cells.Item(7, 9).Value = 0.01857155129685422
$ws.Cells.Item(7, 10).Value = 0.01857155129685422
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 22.04284166666666
$ws.Cells.Item(7, 14).Value = 66.128525
$ws.Cells.Item(7, 15).Value = 0.1704420126740755
$ws.Cells.Item(7, 16).Value = 0.1704420126740755
$ws.Cells.Item(7, 17).Value = 3.087864127261111
$ws.Cells.Item(7, 18).Value = 27.79077714535
$ws.Cells.Item(7, 19).Value = 0.003165372581515669
$ws.Cells.Item(7, 20).Value = 0.00316537258151567

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.778616333333334
$ws.Cells.Item(8, 8).Value = 5.335849000000001
$ws.Cells.Item(8, 9).Value = 0.2357978589514158
$ws.Cells.Item(8, 10).Value = 0.2357978589514158
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 1.321445333333333
$ws.Cells.Item(8, 14).Value = 3.964336
$ws.Cells.Item(8, 15).Value = 0.01021782062667047
$ws.Cells.Item(8, 16).Value = 0.01021782062667047
$ws.Cells.Item(8, 17).Value = 2.350344253473778
$ws.Cells.Item(8, 18).Value = 21.153098281264
$ws.Cells.Item(8, 19).Value = 0.00240934022691851
$ws.Cells.Item(8, 20).Value = 0.00240934022691851

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.778616333333334
$ws.Cells.Item(9, 8).Value = 5.335849000000001
$ws.Cells.Item(9, 9).Value = 0.2357978589514158
$ws.Cells.Item(9, 10).Value = 0.2357978589514158
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 105.9632263333333
$ws.Cells.Item(9, 14).Value = 317.889679
$ws.Cells.Item(9, 15).Value = 0.819340166699254
$ws.Cells.Item(9, 16).Value = 0.8193401666992541
$ws.Cells.Item(9, 17).Value = 188.4679250891635
$ws.Cells.Item(9, 18).Value = 1696.211325802471
$ws.Cells.Item(9, 19).Value = 0.1931986570605802
$ws.Cells.Item(9, 20).Value = 0.1931986570605802

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.778616333333334
$ws.Cells.Item(10, 8).Value = 5.335849000000001
$ws.Cells.Item(10, 9).Value = 0.2357978589514158
$ws.Cells.Item(10, 10).Value = 0.2357978589514158
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 22.04284166666666
$ws.Cells.Item(10, 14).Value = 66.128525
$ws.Cells.Item(10, 15).Value = 0.1704420126740755
$ws.Cells.Item(10, 16).Value = 0.1704420126740755
$ws.Cells.Item(10, 17).Value = 39.20575822141389
$ws.Cells.Item(10, 18).Value = 352.851823992725
$ws.Cells.Item(10, 19).Value = 0.04018986166391707
$ws.Cells.Item(10, 20).Value = 0.04018986166391708
